$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename to short machine-readable codes ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case the Spanish connector words (de/del/el/la/los/y) in state/municipality names ---
$nameChanges = @(
    @('B7', 'Pabellón De Arteaga'),
    @('B8', 'Rincón De Romos'),
    @('B9', 'San José De Gracia'),
    @('B13', 'Playas De Rosarito'),
    @('B40', 'Chiapa De Corzo'),
    @('B45', 'Comitán De Domínguez'),
    @('B65', 'Salto De Agua'),
    @('B98', 'Guadalupe Y Calvo'),
    @('B101', 'Hidalgo Del Parral'),
    @('B120', 'San Francisco De Borja'),
    @('B121', 'San Francisco Del Oro'),
    @('B128', 'Valle De Zaragoza'),
    @('B152', 'Villa De Álvarez'),
    @('A154', 'Ciudad De México'),
    @('B158', 'Cuajimalpa De Morelos'),
    @('B172', 'Coneto De Comonfort'),
    @('B186', 'Nombre De Dios'),
    @('B190', 'Pánuco De Coronado'),
    @('B197', 'San Juan Del Río'),
    @('B198', 'San Pedro Del Gallo'),
    @('A208', 'Estado De México'),
    @('B208', 'Acambay De Ruíz Castañeda'),
    @('B211', 'Almoloya De Juárez'),
    @('B212', 'Almoloya Del Río'),
    @('B218', 'Atizapán De Zaragoza'),
    @('B227', 'Coacalco De Berriozábal'),
    @('B233', 'Ecatepec De Morelos'),
    @('B238', 'Ixtapan De La Sal'),
    @('B239', 'Ixtapan Del Oro'),
    @('B250', 'Naucalpan De Juárez'),
    @('B259', 'San Felipe Del Progreso'),
    @('B260', 'San Simón De Guerrero'),
    @('B270', 'Tenango Del Aire'),
    @('B271', 'Tenango Del Valle'),
    @('B279', 'Tlalnepantla De Baz'),
    @('B284', 'Valle De Bravo'),
    @('B285', 'Valle De Chalco Solidaridad'),
    @('B286', 'Villa De Allende'),
    @('B296', 'San Miguel De Allende'),
    @('B297', 'Apaseo El Alto'),
    @('B298', 'Apaseo El Grande'),
    @('B305', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B309', 'Jaral Del Progreso'),
    @('B317', 'Purísima Del Rincón'),
    @('B321', 'San Diego De La Unión'),
    @('B323', 'San Francisco Del Rincón'),
    @('B325', 'San Luis De La Paz'),
    @('B327', 'Santa Cruz De Juventino Rosas'),
    @('B328', 'Silao De La Victoria'),
    @('B333', 'Valle De Santiago'),
    @('B339', 'Acapulco De Juárez'),
    @('B341', 'Ajuchitlán Del Progreso'),
    @('B342', 'Alcozauca De Guerrero'),
    @('B346', 'Atoyac De Álvarez'),
    @('B347', 'Ayutla De Los Libres'),
    @('B350', 'Buenavista De Cuéllar'),
    @('B351', 'Chilapa De Álvarez'),
    @('B352', 'Chilpancingo De Los Bravo'),
    @('B353', 'Coahuayutla De José María Izazaga'),
    @('B358', 'Coyuca De Benítez'),
    @('B359', 'Coyuca De Catalán'),
    @('B363', 'Cuetzala Del Progreso'),
    @('B364', 'Cutzamala De Pinzón'),
    @('B370', 'Huitzuco De Los Figueroa'),
    @('B371', 'Iguala De La Independencia'),
    @('B373', 'Ixcateopan De Cuauhtémoc'),
    @('B374', 'Zihuatanejo De Azueta'),
    @('B376', 'La Unión De Isidoro Montes De Oca'),
    @('B389', 'Taxco De Alarcón'),
    @('B391', 'Técpan De Galeana'),
    @('B393', 'Tepecoacuilco De Trujano'),
    @('B395', 'Tixtla De Guerrero'),
    @('B399', 'Tlalixtaquilla De Maldonado'),
    @('B400', 'Tlapa De Comonfort'),
    @('B410', 'Agua Blanca De Iturbide'),
    @('B416', 'Atotonilco De Tula'),
    @('B417', 'Atotonilco El Grande'),
    @('B422', 'Cuautepec De Hinojosa'),
    @('B426', 'Huasca De Ocampo'),
    @('B428', 'Huejutla De Reyes'),
    @('B437', 'Mineral De La Reforma'),
    @('B438', 'Mineral Del Monte'),
    @('B439', 'Mixquiahuala De Juárez'),
    @('B440', 'Molango De Escamilla'),
    @('B442', 'Nopala De Villagrán'),
    @('B443', 'Pachuca De Soto'),
    @('B444', 'Progreso De Obregón'),
    @('B449', 'Santiago De Anaya'),
    @('B452', 'Tenango De Doria'),
    @('B454', 'Tepeji Del Río De Ocampo'),
    @('B457', 'Tezontepec De Aldama'),
    @('B465', 'Tula De Allende'),
    @('B466', 'Tulancingo De Bravo'),
    @('B467', 'Villa De Tezontepec'),
    @('B471', 'Zacualtipán De Ángeles'),
    @('B472', 'Zapotlán De Juárez'),
    @('B477', 'Acatlán De Juárez'),
    @('B478', 'Ahualulco De Mercado'),
    @('B483', 'Atemajac De Brizuela'),
    @('B486', 'Atotonilco El Alto'),
    @('B488', 'Autlán De Navarro'),
    @('B494', 'Cañadas De Obregón'),
    @('B501', 'Concepción De Buenos Aires'),
    @('B502', 'Cuautitlán De García Barragán'),
    @('B510', 'Encarnación De Díaz'),
    @('B517', 'Huejuquilla El Alto'),
    @('B518', 'Ixtlahuacán De Los Membrillos'),
    @('B519', 'Ixtlahuacán Del Río'),
    @('B523', 'Jilotlán De Los Dolores'),
    @('B528', 'La Manzanilla De La Paz'),
    @('B529', 'Lagos De Moreno'),
    @('B537', 'Ojuelos De Jalisco'),
    @('B542', 'San Cristóbal De La Barranca'),
    @('B543', 'San Diego De Alejandría'),
    @('B545', 'San Juan De Los Lagos'),
    @('B546', 'San Juanito De Escobedo'),
    @('B549', 'San Martín De Bolaños'),
    @('B551', 'San Miguel El Alto'),
    @('B552', 'San Sebastián Del Oeste'),
    @('B553', 'Santa María De Los Ángeles'),
    @('B554', 'Santa María Del Oro'),
    @('B557', 'Talpa De Allende'),
    @('B558', 'Tamazula De Gordiano'),
    @('B564', 'Teocuitatlán De Corona'),
    @('B565', 'Tepatitlán De Morelos'),
    @('B568', 'Tizapán El Alto'),
    @('B569', 'Tlajomulco De Zúñiga'),
    @('B581', 'Unión De San Antonio'),
    @('B582', 'Unión De Tula'),
    @('B583', 'Valle De Juárez'),
    @('B588', 'Yahualica De González Gallo'),
    @('B589', 'Zacoalco De Torres'),
    @('B592', 'Zapotitlán De Vadillo'),
    @('B593', 'Zapotlán Del Rey'),
    @('B594', 'Zapotlán El Grande'),
    @('B619', 'Coalcomán De Vázquez Pallares'),
    @('B621', 'Cojumatlán De Régules'),
    @('B688', 'Tiquicheo De Nicolás Romero'),
    @('B714', 'Coatlán Del Río'),
    @('B721', 'Jonacatepec De Leandro Valle'),
    @('B724', 'Puente De Ixtla'),
    @('B730', 'Tlaltizapán De Zapata'),
    @('B737', 'Zacualpan De Amilpas'),
    @('B741', 'Amatlán De Cañas'),
    @('B742', 'Bahía De Banderas'),
    @('B746', 'Ixtlán Del Río'),
    @('B753', 'Santa María Del Oro'),
    @('B771', 'San Nicolás De Los Garza'),
    @('B774', 'Acatlán De Pérez Figueroa'),
    @('B776', 'Ayoquezco De Aldama'),
    @('B779', 'Chalcatongo De Hidalgo'),
    @('B782', 'Constancia Del Rosario'),
    @('B784', 'Cuyamecalco Villa De Zaragoza'),
    @('B785', 'Guevea De Humboldt'),
    @('B786', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B787', 'Heroica Ciudad De Huajuapan De León'),
    @('B788', 'Heroica Ciudad De Tlaxiaco'),
    @('B789', 'Ixtlán De Juárez'),
    @('B790', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B794', 'Mariscala De Juárez'),
    @('B796', 'Miahuatlán De Porfirio Díaz'),
    @('B797', 'Nejapa De Madero'),
    @('B798', 'Oaxaca De Juárez'),
    @('B799', 'Ocotlán De Morelos'),
    @('B800', 'Putla Villa De Guerrero'),
    @('B801', 'Reforma De Pineda'),
    @('B802', 'Rojas De Cuauhtémoc'),
    @('B808', 'San Antonio De La Cal'),
    @('B824', 'San José Del Progreso'),
    @('B832', 'San Juan Del Estado'),
    @('B880', 'Santa Cruz Tacache De Mina'),
    @('B917', 'Tamazulápam Del Espíritu Santo'),
    @('B918', 'Tataltepec De Valdés'),
    @('B919', 'Teococuilco De Marcos Pérez'),
    @('B920', 'Teotitlán De Flores Magón'),
    @('B921', 'Teotitlán Del Valle'),
    @('B923', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B924', 'Tlacolula De Matamoros'),
    @('B925', 'Totontepec Villa De Morelos'),
    @('B926', 'Villa De Etla'),
    @('B927', 'Villa De Tututepec'),
    @('B928', 'Villa De Zaachila'),
    @('B930', 'Villa Sola De Vega'),
    @('B931', 'Zapotitlán Del Río'),
    @('B933', 'Zimatlán De Álvarez'),
    @('B954', 'Ayotoxco De Guerrero'),
    @('B956', 'Chalchicomula De Sesma'),
    @('B965', 'Chila De La Sal'),
    @('B973', 'Cuapiaxtla De Madero'),
    @('B976', 'Cuayuca De Andrade'),
    @('B977', 'Cuetzalan Del Progreso'),
    @('B985', 'Huehuetlán El Chico'),
    @('B986', 'Huehuetlán El Grande'),
    @('B989', 'Ixcamilpa De Guerrero'),
    @('B993', 'Izúcar De Matamoros'),
    @('B1000', 'Los Reyes De Juárez'),
    @('B1007', 'Palmar De Bravo'),
    @('B1022', 'San Nicolás De Los Ranchos'),
    @('B1026', 'San Salvador El Seco'),
    @('B1027', 'San Salvador El Verde'),
    @('B1033', 'Tecali De Herrera'),
    @('B1039', 'Tepanco De López'),
    @('B1040', 'Tepatlaxco De Hidalgo'),
    @('B1044', 'Tepexi De Rodríguez'),
    @('B1046', 'Tetela De Ocampo'),
    @('B1051', 'Tlacotepec De Benito Juárez'),
    @('B1066', 'Xochitlán De Vicente Suárez'),
    @('B1078', 'Amealco De Bonfil'),
    @('B1079', 'Cadereyta De Montes'),
    @('B1085', 'Pinal De Amoles'),
    @('B1088', 'San Juan Del Río'),
    @('B1100', 'Axtla De Terrazas'),
    @('B1109', 'Mexquitic De Carmona'),
    @('B1112', 'San Ciro De Acosta'),
    @('B1116', 'Santa María Del Río'),
    @('B1124', 'Villa De Arriaga'),
    @('B1125', 'Villa De Ramos'),
    @('B1169', 'Nacozari De García'),
    @('B1190', 'Jalpa De Méndez'),
    @('B1208', 'Acuamanala De Miguel Hidalgo'),
    @('B1212', 'Contla De Juan Cuamatzi'),
    @('B1215', 'Muñoz De Domingo Arenas'),
    @('B1216', 'Nanacamilpa De Mariano Arista'),
    @('B1219', 'Papalotla De Xicohténcatl'),
    @('B1223', 'Tetla De La Solidaridad'),
    @('B1238', 'Amatlán De Los Reyes'),
    @('B1245', 'Boca Del Río'),
    @('B1248', 'Castillo De Teayo'),
    @('B1257', 'Cosamaloapan De Carpio'),
    @('B1258', 'Cosautlán De Carvajal'),
    @('B1270', 'Hueyapan De Ocampo'),
    @('B1271', 'Huiloapan De Cuauhtémoc'),
    @('B1272', 'Ignacio De La Llave'),
    @('B1274', 'Ixhuatlán De Madero'),
    @('B1275', 'Ixhuatlán Del Café'),
    @('B1276', 'Ixhuatlán Del Sureste'),
    @('B1284', 'Juchique De Ferrer'),
    @('B1287', 'Lerdo De Tejada'),
    @('B1292', 'Martínez De La Torre'),
    @('B1293', 'Medellín De Bravo'),
    @('B1304', 'Paso De Ovejas'),
    @('B1308', 'Poza Rica De Hidalgo'),
    @('B1316', 'Sayula De Alemán'),
    @('B1319', 'Soledad De Doblado'),
    @('B1339', 'Vega De Alatorre'),
    @('B1347', 'Zozocolco De Hidalgo'),
    @('B1353', 'Dzilam De Bravo'),
    @('B1372', 'Concepción Del Oro'),
    @('B1374', 'El Plateado De Joaquín Amaro'),
    @('B1384', 'Jiménez Del Teul'),
    @('B1390', 'Mezquital Del Oro'),
    @('B1395', 'Moyahua De Estrada'),
    @('B1396', 'Nochistlán De Mejía'),
    @('B1397', 'Noria De Ángeles'),
    @('B1408', 'Teúl De González Ortega'),
    @('B1409', 'Tlaltenango De Sánchez Román'),
    @('B1410', 'Trinidad García De La Cadena'),
    @('B1413', 'Villa De Cos')
)
foreach ($pair in $nameChanges) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Tiny floating point precision refresh on recomputed percentage cells ---
$pctChanges = @(
    @('D10', 0.009935644123292313),
    @('D154', 0.009540476459297733),
    @('D202', 0.0009032403748447556),
    @('D343', 0.0009032403748447556),
    @('D358', 0.0009032403748447556),
    @('D383', 0.0009032403748447556),
    @('D421', 0.009314666365586539),
    @('D479', 0.0009032403748447556),
    @('D495', 0.0009032403748447556),
    @('D499', 0.0009032403748447556),
    @('D502', 0.0009032403748447556),
    @('D619', 0.0009032403748447556),
    @('D652', 0.009596928982725529),
    @('D933', 0.0009032403748447556),
    @('D1017', 0.009484023935869931),
    @('D1060', 0.0009032403748447556),
    @('D1086', 0.0009032403748447556),
    @('D1132', 0.0009032403748447556),
    @('D1157', 0.0009032403748447556),
    @('D1359', 0.0009032403748447556),
    @('D1406', 0.0009032403748447556)
)
foreach ($pair in $pctChanges) {
    $ws.Range($pair[0]).Value2 = $pair[1]
}

# --- Drop the trailing footer/metadata rows (1422-1426); dimension shrinks to A1:D1420 ---
$ws.Range("A1422:A1426").EntireRow.Delete() | Out-Null

Write-Host "edit complete"
